$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add "Yes" markers in column D for the newly completed items (d46 calc + background)
$ws.Range("D7").Value = "Yes"
$ws.Range("D8").Value = "Yes"
$ws.Range("D22").Value = "Yes"
$ws.Range("D23").Value = "Yes"

# Update the active selection to reflect where the user left off
$ws.Range("D18").Select()
